$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the 1-based index of the first paragraph whose text contains
# the given substring.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($doc, $needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# Change 1: cover-page "Version 1.0" run picks up redundant theme tint/shade
# attributes on its color, and the following paragraph that reads
# "{Name} (Project Manager)" is removed entirely.
# ---------------------------------------------------------------------------
$verIdx = Find-ParagraphIndex $d "Version 1.0"
if ($verIdx -ge 1) {
    $pVer = $d.Paragraphs.Item($verIdx)
    $rVer = $d.Range($pVer.Range.Start, $pVer.Range.End)
    $verXml = $pkgHeader + '<w:p><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="48"/><w:szCs w:val="48"/><w:vertAlign w:val="superscript"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="48"/><w:szCs w:val="48"/><w:vertAlign w:val="superscript"/><w:lang w:val="en-US"/></w:rPr><w:t>Version 1.0</w:t></w:r></w:p>' + $pkgFooter
    $rVer.InsertXML($verXml)

    # The very next paragraph is "{Name} (Project Manager)" - drop it completely.
    $nameIdx = $verIdx + 1
    $pName = $d.Paragraphs.Item($nameIdx)
    if ($pName.Range.Text.Contains("{Name}")) {
        $rName = $d.Range($pName.Range.Start, $pName.Range.End)
        $rName.Delete()
    }
}

# ---------------------------------------------------------------------------
# Change 2: the empty, underlined/italic "note" placeholder paragraph right
# after the "**Note: Please review..." text gains a pStyle + explicit
# i/iCs val + theme tint/shade on its color, and the following two
# paragraphs (the "Note: ... ReactJS ... Node.js ... Java server ..."
# paragraph and the blank paragraph immediately after it) are removed.
# ---------------------------------------------------------------------------
$reactIdx = Find-ParagraphIndex $d "Note: The UI and backend are being updated"
if ($reactIdx -ge 2) {
    $placeholderIdx = $reactIdx - 1
    $pPlaceholder = $d.Paragraphs.Item($placeholderIdx)
    $rPlaceholder = $d.Range($pPlaceholder.Range.Start, $pPlaceholder.Range.End)
    $placeholderXml = $pkgHeader + '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:after="0" w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:eastAsia="Arial" w:cs="Arial"/><w:i w:val="1"/><w:iCs w:val="1"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' + $pkgFooter
    $rPlaceholder.InsertXML($placeholderXml)

    # Delete the "Note: ... ReactJS ..." paragraph and the blank paragraph
    # directly following it (two paragraphs total).
    $pReact = $d.Paragraphs.Item($reactIdx)
    $pAfterReact = $d.Paragraphs.Item($reactIdx + 1)
    $rDelete = $d.Range($pReact.Range.Start, $pAfterReact.Range.End)
    $rDelete.Delete()
}
